$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 322, pushing old rows 322-324 down to 326-328.
$ws.Range("A322:T325").Insert()

# --- Row 322: updated "Carson / Especial" record (was Doctor Davis / Especial) ---
$ws.Cells.Item(322, 1).Value = 7
$ws.Cells.Item(322, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(322, 3).Value = "Ñuble"
$ws.Cells.Item(322, 4).Value = 44939
$ws.Cells.Item(322, 5).Value = 16
$ws.Cells.Item(322, 6).Value = "Fruta"
$ws.Cells.Item(322, 7).Value = 100103
$ws.Cells.Item(322, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(322, 9).Value = 100103004
$ws.Cells.Item(322, 10).Value = "Durazno"
$ws.Cells.Item(322, 11).Value = "Carson"
$ws.Cells.Item(322, 12).Value = "Especial"
$ws.Cells.Item(322, 13).Value = 60
$ws.Cells.Item(322, 14).Value = 18000
$ws.Cells.Item(322, 15).Value = 18000
$ws.Cells.Item(322, 16).Value = 18000
$ws.Cells.Item(322, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(322, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(322, 19).Value = 1200
$ws.Cells.Item(322, 20).Value = 15

# --- Row 323: updated "Carson / Primera" record (was Doctor Davis / Primera) ---
$ws.Cells.Item(323, 1).Value = 7
$ws.Cells.Item(323, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(323, 3).Value = "Ñuble"
$ws.Cells.Item(323, 4).Value = 44939
$ws.Cells.Item(323, 5).Value = 16
$ws.Cells.Item(323, 6).Value = "Fruta"
$ws.Cells.Item(323, 7).Value = 100103
$ws.Cells.Item(323, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(323, 9).Value = 100103004
$ws.Cells.Item(323, 10).Value = "Durazno"
$ws.Cells.Item(323, 11).Value = "Carson"
$ws.Cells.Item(323, 12).Value = "Primera"
$ws.Cells.Item(323, 13).Value = 120
$ws.Cells.Item(323, 14).Value = 16000
$ws.Cells.Item(323, 15).Value = 17000
$ws.Cells.Item(323, 16).Value = 16500
$ws.Cells.Item(323, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(323, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(323, 19).Value = 1100
$ws.Cells.Item(323, 20).Value = 15

# --- Row 324: new "Carson / Segunda" record ---
$ws.Cells.Item(324, 1).Value = 7
$ws.Cells.Item(324, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(324, 3).Value = "Ñuble"
$ws.Cells.Item(324, 4).Value = 44939
$ws.Cells.Item(324, 5).Value = 16
$ws.Cells.Item(324, 6).Value = "Fruta"
$ws.Cells.Item(324, 7).Value = 100103
$ws.Cells.Item(324, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(324, 9).Value = 100103004
$ws.Cells.Item(324, 10).Value = "Durazno"
$ws.Cells.Item(324, 11).Value = "Carson"
$ws.Cells.Item(324, 12).Value = "Segunda"
$ws.Cells.Item(324, 13).Value = 120
$ws.Cells.Item(324, 14).Value = 14000
$ws.Cells.Item(324, 15).Value = 14000
$ws.Cells.Item(324, 16).Value = 14000
$ws.Cells.Item(324, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(324, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(324, 19).Value = 933
$ws.Cells.Item(324, 20).Value = 15

# --- Row 325: new "Kurakata / Especial" record ---
$ws.Cells.Item(325, 1).Value = 7
$ws.Cells.Item(325, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(325, 3).Value = "Ñuble"
$ws.Cells.Item(325, 4).Value = 44939
$ws.Cells.Item(325, 5).Value = 16
$ws.Cells.Item(325, 6).Value = "Fruta"
$ws.Cells.Item(325, 7).Value = 100103
$ws.Cells.Item(325, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(325, 9).Value = 100103004
$ws.Cells.Item(325, 10).Value = "Durazno"
$ws.Cells.Item(325, 11).Value = "Kurakata"
$ws.Cells.Item(325, 12).Value = "Especial"
$ws.Cells.Item(325, 13).Value = 80
$ws.Cells.Item(325, 14).Value = 18000
$ws.Cells.Item(325, 15).Value = 18000
$ws.Cells.Item(325, 16).Value = 18000
$ws.Cells.Item(325, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(325, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(325, 19).Value = 1200
$ws.Cells.Item(325, 20).Value = 15
